# "Generate Report for Handoff"
#
# The localization-status workbook is refreshed with a new handoff report:
#  - the previously-failing source file was renamed and has now transformed
#    successfully ("Ready for handoff"), with real handoff target files and
#    timestamps recorded for both locales;
#  - a brand-new source file shows up with the very same "ready" status;
#  - the untouched .localization-config bookkeeping row simply moves down to
#    make room for the new entry.

$wb = $excel.ActiveWorkbook

$oldFile = "d473859f-d9d0-42d3-96c3-ffddf165fab4.md"
$renamedFile = "fe289fc4-9e55-4194-a23f-dd16b44051a8.md"
$newFile = "ffff6b28143b-acf1-4a7d-8773-f644510d5524.md"
$configFile = ".localization-config"

$readyStatus = "Ready for handoff"
$notLocalizedStatus = "Not to be localized"

$zhHandoffFile = "fe289fc4-9e55-4194-a23f-dd16b44051a8.0cec61a9d3e388cacdf259b410cc67ac11c93aed.zh-cn.xlf"
$deHandoffFile = "fe289fc4-9e55-4194-a23f-dd16b44051a8.0cec61a9d3e388cacdf259b410cc67ac11c93aed.de-de.xlf"
$zhHandoffDate = "2016-02-17 06:34:47"
$deHandoffDate = "2016-02-17 06:34:58"
$epoch = "0001-01-01 00:00:00"
$includeReason = "Include"
$ignoredReason = "Ignored"

$renamedUrl = "https://github.com/OpenLocalizationTest/oltest/blob/fe289fc49e554194a23fdd16b44051a8/e2e/$renamedFile"
$newUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ffff6b28143bacf14a7d8773f644510d5524/e2e/$newFile"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d4180c1b81328e9e85e00746173f61cec526913b/$configFile"
$zhHandoffUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0cec61a9d3e388cacdf259b410cc67ac11c93aed/e2e/$zhHandoffFile"
$deHandoffUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0cec61a9d3e388cacdf259b410cc67ac11c93aed/e2e/$deHandoffFile"

# ---------------------------------------------------------------------------
# Overview sheet: File Name / zh-cn / de-de summary
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows.Item(3).Insert()

$wsOverview.Range("A2").Value = $renamedFile
$wsOverview.Range("B2").Value = $readyStatus
$wsOverview.Range("C2").Value = $readyStatus

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus

# Row 4 keeps the untouched ".localization-config" values that Insert()
# already carried down from the old row 3.

$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $renamedUrl, "", "", $renamedFile) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $newUrl, "", "", $newFile) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $configUrl, "", "", $configFile) | Out-Null

# ---------------------------------------------------------------------------
# Per-locale detail sheets: zh-cn and de-de
# ---------------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; HandoffFile = $zhHandoffFile; HandoffDate = $zhHandoffDate; HandoffUrl = $zhHandoffUrl },
    @{ Sheet = "de-de"; HandoffFile = $deHandoffFile; HandoffDate = $deHandoffDate; HandoffUrl = $deHandoffUrl }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    $ws.Rows.Item(3).Insert()

    # Row 2: renamed source file, transform now succeeds
    $ws.Range("A2").Value = $renamedFile
    $ws.Range("B2").Value = $readyStatus
    $ws.Range("C2").Value = $locale.HandoffFile
    $ws.Range("D2").Value = $locale.HandoffDate
    $ws.Range("G2").Value = $epoch
    $ws.Range("H2").Value = $includeReason

    # Row 3 (new): additional source file, same handoff batch
    $ws.Range("A3").Value = $newFile
    $ws.Range("B3").Value = $readyStatus
    $ws.Range("C3").Value = $locale.HandoffFile
    $ws.Range("D3").Value = $locale.HandoffDate
    $ws.Range("G3").Value = $epoch
    $ws.Range("H3").Value = $includeReason

    # Row 4 keeps the untouched ".localization-config" bookkeeping values.

    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $renamedUrl, "", "", $renamedFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $locale.HandoffUrl, "", "", $locale.HandoffFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $newUrl, "", "", $newFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $locale.HandoffUrl, "", "", $locale.HandoffFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $configFile) | Out-Null
}
